# Update workbook with refreshed TPM-derived NATMI values.
# This mirrors a re-run of the upstream analysis script which produced
# slightly different ligand/receptor expression, specificity, and edge
# weight statistics for the Bmp4-Bmpr1b sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> FAPs)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.202518666666666
$ws.Range("H2").Value = 12.607556
$ws.Range("I2").Value = 0.08075097102331126
$ws.Range("J2").Value = 0.08075097102331129
$ws.Range("Q2").Value = 5.432748571911555
$ws.Range("R2").Value = 48.894737147204
$ws.Range("S2").Value = 0.07287201082608399
$ws.Range("T2").Value = 0.07287201082608401

# Row 3 (ECs -> MuSCs)
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.202518666666666
$ws.Range("H3").Value = 12.607556
$ws.Range("I3").Value = 0.08075097102331126
$ws.Range("J3").Value = 0.08075097102331129
$ws.Range("Q3").Value = 0.5873916373982222
$ws.Range("R3").Value = 5.286524736584
$ws.Range("S3").Value = 0.007878960197227274
$ws.Range("T3").Value = 0.007878960197227275

# Row 4 (FAPs -> FAPs)
$ws.Range("I4").Value = 0.7496282157262072
$ws.Range("J4").Value = 0.7496282157262073
$ws.Range("S4").Value = 0.6764861742178735
$ws.Range("T4").Value = 0.6764861742178736

# Row 5 (FAPs -> MuSCs)
$ws.Range("I5").Value = 0.7496282157262072
$ws.Range("J5").Value = 0.7496282157262073
$ws.Range("S5").Value = 0.07314204150833373
$ws.Range("T5").Value = 0.07314204150833374

# Row 6 (MuSCs -> FAPs)
$ws.Range("G6").Value = 8.827567333333333
$ws.Range("I6").Value = 0.1696208132504815
$ws.Range("J6").Value = 0.1696208132504815
$ws.Range("S6").Value = 0.1530707257495391
$ws.Range("T6").Value = 0.1530707257495391

# Row 7 (MuSCs -> MuSCs)
$ws.Range("G7").Value = 8.827567333333333
$ws.Range("I7").Value = 0.1696208132504815
$ws.Range("J7").Value = 0.1696208132504815
$ws.Range("S7").Value = 0.01655008750094238
$ws.Range("T7").Value = 0.01655008750094239

$wb.Save()
